# Generate Report for Handback
#
# For the "7d38ce2a-0449-4b21-88aa-c427a8c74d57" entry (row 8) in both the
# zh-cn and de-de localization-status worksheets, the handback-generation
# job filled in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns and recorded that the handed-back
# file version was stale in the "Error Detail" column.

$wb = $excel.ActiveWorkbook

$latestHandbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2f5c7af7cd1c55114c8f4d36cecbd1e774728f1e/e2e/7d38ce2a-0449-4b21-88aa-c427a8c74d57.md"
$mdDisplay = "7d38ce2a-0449-4b21-88aa-c427a8c74d57.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96a07c92a241e178fed8e48676350a0b0365c0ee/e2e/7d38ce2a-0449-4b21-88aa-c427a8c74d57.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2f5c7af7cd1c55114c8f4d36cecbd1e774728f1e/e2e/7d38ce2a-0449-4b21-88aa-c427a8c74d57.md."

function Update-LocaleSheet($SheetName, $TargetXlf, $HandbackDateTime) {
    $ws = $wb.Worksheets.Item($SheetName)

    # I8 - Latest Target File: becomes a hyperlink pointing at the (latest)
    # handback markdown file, same as column A's hyperlink for this row.
    $ws.Hyperlinks.Add($ws.Range("I8"), $latestHandbackUrl, "", "", $mdDisplay) | Out-Null

    # J8 - Latest Handback File.
    $ws.Range("J8").Value = $TargetXlf

    # K8 - Latest Handback DateTime.
    $ws.Range("K8").Value = $HandbackDateTime

    # P8 - Error Detail: the handback version-mismatch warning.
    $ws.Range("P8").Value = $errorDetail

    # Error Detail column (P / 16) is widened so the message is readable.
    $ws.Columns.Item(16).ColumnWidth = 39.17
}

Update-LocaleSheet "zh-cn" "7d38ce2a-0449-4b21-88aa-c427a8c74d57.3e7fcb012ce02ef9aabd9015e14b7133c2a40e6f.zh-cn.xlf" "2016-09-04 06:47:08"
Update-LocaleSheet "de-de" "7d38ce2a-0449-4b21-88aa-c427a8c74d57.3e7fcb012ce02ef9aabd9015e14b7133c2a40e6f.de-de.xlf" "2016-09-04 06:47:16"
